$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 182
$ws.Range("D21").Value = 162
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = 46.41833810888252
